{"js": "// Update the 25 \"a\u00f7b=\" arithmetic problems in the practice-sheet table to a\n// new set of values (same pattern, new operands), per the commit's\n// regenerated worksheet.\nconst replacements = [\n  [\"676\u00f78=\", \"792\u00f72=\"],\n  [\"477\u00f74=\", \"189\u00f77=\"],\n  [\"208\u00f75=\", \"691\u00f77=\"],\n  [\"910\u00f72=\", \"449\u00f73=\"],\n  [\"308\u00f73=\", \"904\u00f72=\"],\n  [\"754\u00f75=\", \"594\u00f72=\"],\n  [\"636\u00f78=\", \"472\u00f77=\"],\n  [\"654\u00f79=\", \"441\u00f75=\"],\n  [\"259\u00f72=\", \"337\u00f73=\"],\n  [\"874\u00f76=\", \"470\u00f72=\"],\n  [\"209\u00f76=\", \"271\u00f73=\"],\n  [\"676\u00f74=\", \"503\u00f77=\"],\n  [\"107\u00f73=\", \"335\u00f78=\"],\n  [\"834\u00f77=\", \"582\u00f76=\"],\n  [\"170\u00f72=\", \"285\u00f77=\"],\n  [\"542\u00f78=\", \"630\u00f78=\"],\n  [\"626\u00f75=\", \"838\u00f78=\"],\n  [\"591\u00f76=\", \"927\u00f79=\"],\n  [\"758\u00f77=\", \"296\u00f79=\"],\n  [\"582\u00f75=\", \"782\u00f77=\"],\n  [\"483\u00f76=\", \"769\u00f73=\"],\n  [\"261\u00f77=\", \"455\u00f74=\"],\n  [\"454\u00f78=\", \"986\u00f79=\"],\n  [\"410\u00f75=\", \"545\u00f75=\"],\n  [\"750\u00f74=\", \"251\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 \"a\u00f7b=\" arithmetic problems in the practice-sheet table to a\n# new set of values (same pattern, new operands), per the commit's\n# regenerated worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"676\u00f78=\", \"792\u00f72=\"),\n  @(\"477\u00f74=\", \"189\u00f77=\"),\n  @(\"208\u00f75=\", \"691\u00f77=\"),\n  @(\"910\u00f72=\", \"449\u00f73=\"),\n  @(\"308\u00f73=\", \"904\u00f72=\"),\n  @(\"754\u00f75=\", \"594\u00f72=\"),\n  @(\"636\u00f78=\", \"472\u00f77=\"),\n  @(\"654\u00f79=\", \"441\u00f75=\"),\n  @(\"259\u00f72=\", \"337\u00f73=\"),\n  @(\"874\u00f76=\", \"470\u00f72=\"),\n  @(\"209\u00f76=\", \"271\u00f73=\"),\n  @(\"676\u00f74=\", \"503\u00f77=\"),\n  @(\"107\u00f73=\", \"335\u00f78=\"),\n  @(\"834\u00f77=\", \"582\u00f76=\"),\n  @(\"170\u00f72=\", \"285\u00f77=\"),\n  @(\"542\u00f78=\", \"630\u00f78=\"),\n  @(\"626\u00f75=\", \"838\u00f78=\"),\n  @(\"591\u00f76=\", \"927\u00f79=\"),\n  @(\"758\u00f77=\", \"296\u00f79=\"),\n  @(\"582\u00f75=\", \"782\u00f77=\"),\n  @(\"483\u00f76=\", \"769\u00f73=\"),\n  @(\"261\u00f77=\", \"455\u00f74=\"),\n  @(\"454\u00f78=\", \"986\u00f79=\"),\n  @(\"410\u00f75=\", \"545\u00f75=\"),\n  @(\"750\u00f74=\", \"251\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
